# produits.xlsx update:
#  - rename sheet "Produits" -> "Sheet1"
#  - lowercase/rename a couple of header labels and add a new "mode_d_application" column
#  - swap rows 2/3 (ortiva now comes before vertimec) and fill in the new column
#  - drop the trailing "Test/Hiba/Zwina" scratch row
#  - bold + center + box-border the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Sheet1"

# Header row (A1 keeps its text; B1/C1 get new casing, D1 is brand new)
$ws.Range("B1").Value = "dose"
$ws.Range("C1").Value = "cible"
$ws.Range("D1").Value = "mode_d_application"

# Row 2 becomes the old "ortiva" product line, plus the new column
$ws.Range("A2").Value = "ortiva"
$ws.Range("B2").Value = "50cc/hl"
$ws.Range("C2").Value = "oidium"
$ws.Range("D2").Value = "feuilles"

# Row 3 becomes the old "vertimec" product line, plus the new column
$ws.Range("A3").Value = "vertimec"
$ws.Range("B3").Value = "50cc/hl"
$ws.Range("C3").Value = "oidium"
$ws.Range("D3").Value = "feuilles"

# Remove the old trailing test row entirely
$ws.Rows("4").Delete()

# Style the header row: bold, centered, thin box border
$header = $ws.Range("A1:D1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1
